$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("cs_mom_lo")
$rowVals = @(2.678779653636064, 2.948158636367531, 3.039929902928181, 3.581280847811459, 2.775792493577502, 2.993241945912671, 2.978520547398343, 3.207118279911883, 2.883352323197843, 2.916563500049703, 3.270631217356445, 3.768162223565663, 2.856264157795118, 2.726541848389753, 3.045078341994438, 3.66555241180345, 2.983325283614649, 2.88094268993568, 3.36869051898925, 3.814176548268978, 2.906732973915696, 2.789272129948952, 3.09054018779733, 3.637670385609148)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.09232732632062257, 0.09757329529790892, 0.09928446386480916, 0.1086964891945552, 0.09425664977233983, 0.09841850023176724, 0.098143491391951, 0.1023107656695093, 0.09634243278747578, 0.09697558269458306, 0.1034311122053919, 0.1117058726205788, 0.09582226545224026, 0.09328303839456398, 0.0993793847946034, 0.1100674311620693, 0.09823335227337138, 0.09629629848105692, 0.1051307064551177, 0.112429965276351, 0.09678869309198257, 0.09452105840398439, 0.1002127034723048, 0.1096164158310571)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.1456915117055111, 0.1399893340782122, 0.1393993029592432, 0.1346396932034567, 0.1447300651247116, 0.1384866236150863, 0.1394439803551607, 0.1349431124270856, 0.1440964800874713, 0.1416718409430466, 0.1391413615539735, 0.1348472207937331, 0.1432955180510785, 0.1396686668283787, 0.1383534989431277, 0.1345420060764121, 0.1434024767262408, 0.1410224553229236, 0.1387354020214656, 0.1358491468349574, 0.1425784018376048, 0.1395473161456954, 0.1375968673350875, 0.1349500731384396)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.633717951305533, 0.6970052107212297, 0.7122307052987015, 0.8073138508292783, 0.6512582557820339, 0.7106715266978738, 0.7038202089612023, 0.7581770112556966, 0.6685967119321221, 0.6845085237056258, 0.7433527388997883, 0.8283883936432622, 0.6687038558881088, 0.6678880847999201, 0.718300480679963, 0.8180897131826422, 0.6850185193168004, 0.6828437234378777, 0.7577785116365, 0.8276089169182761, 0.6788454060680511, 0.6773405681647003, 0.7283065771276491, 0.8122738527055572)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(-0.2520100508991268, -0.2128140678433241, -0.2250117050171943, -0.167542326387377, -0.2516989295527468, -0.2110967793378163, -0.2153573872247184, -0.1745615337730845, -0.2319106643025202, -0.2151451793192236, -0.2044095511908439, -0.1686515657583613, -0.2349010206482688, -0.2183405212429793, -0.2110825267419936, -0.1708265750765352, -0.2294223355660262, -0.2183344143304764, -0.2049245767052084, -0.1759903235940466, -0.2328208971208878, -0.2188207465453388, -0.2049614630743272, -0.1767902823749796)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value2 = $rowVals[$i] }

$ws = $wb.Worksheets.Item("cs_idio_mom_lo")
$rowVals = @(2.827512169857502, 3.076093117959826, 3.020725331876215, 3.600534620734132, 3.036618406120152, 3.055892359113876, 3.185467326986447, 3.650614002798783, 2.868704075646829, 3.079346815707695, 3.041150118655348, 3.625430572674657, 3.005664119016972, 3.160529927676158, 3.075140276535982, 3.496175713262221, 2.932772387561908, 3.086036780573321, 3.08670304456857, 3.485612264727372, 2.958814747705048, 3.054796803204744, 3.079642251240935, 3.51761327336855)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.09526640840684419, 0.09994882758040791, 0.09892939466262574, 0.109011772019221, 0.09922335072975752, 0.09957839489495357, 0.1019252444672385, 0.1098261032226231, 0.09606156765541063, 0.1000083323858769, 0.09930697097440211, 0.1094176276411234, 0.09864982326057192, 0.1014788936349946, 0.09993139328832612, 0.1072879175882233, 0.09728277946868902, 0.1001305420121645, 0.1001427028436261, 0.1071113516078219, 0.09777388097744888, 0.09955825593738288, 0.1000137332269018, 0.1076450566484293)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.1462098117183996, 0.1425741882184576, 0.1418467175436122, 0.1383885570825256, 0.1458462171869417, 0.1425772548193096, 0.1410935285757994, 0.1378922868301222, 0.1439056593809458, 0.1407895267381752, 0.1398471991326542, 0.1378234148324516, 0.1438433733372125, 0.1405556585279092, 0.1395032936549105, 0.1370669179516612, 0.1435273375356112, 0.1405469405122318, 0.1393696400018997, 0.1383699960105145, 0.1428850095132435, 0.1401318378722584, 0.1393045625445846, 0.1378293477893488)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.6515732924294265, 0.7010303115123651, 0.6974387308765807, 0.7877224412002046, 0.6803285861201028, 0.6984171144349134, 0.722394892920134, 0.7964629911310738, 0.6675315485759827, 0.7103392894548282, 0.710110546298485, 0.7938972327317504, 0.6858141669780421, 0.7219836945578726, 0.7163371607234345, 0.7827411544050336, 0.6777996522408266, 0.7124348751188234, 0.7185402993238779, 0.7740937681293473, 0.6842836859550796, 0.710461358739463, 0.7179501618612973, 0.7810024379782191)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(-0.2438900034816637, -0.2257452471135075, -0.20406299126646, -0.1881908161494522, -0.2334214482113403, -0.2048050215652706, -0.1964748987888418, -0.1993829151180535, -0.2188999678137772, -0.2055052968019089, -0.1975081114184588, -0.1911585025624404, -0.2162717979524342, -0.196038836500656, -0.1954966515902944, -0.1947289925060446, -0.2215535525923328, -0.2083442298189825, -0.1948898637228887, -0.1881098454521793, -0.2146228959917894, -0.2028443198050048, -0.1950109338199709, -0.1968102857935684)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value2 = $rowVals[$i] }

$ws = $wb.Worksheets.Item("cs_reversal_lo")
$rowVals = @(2.77499742049895, 3.066854469023514, 3.555879891793272, 2.892607977107825, 2.746935875089095, 2.840234129872154, 3.390542867058371, 3.014466727917439, 2.740420346005984, 2.78790743885687, 3.031122247959569, 2.837146812274931, 2.765003528729698, 2.730046539845006, 3.016695152645311, 2.948930029484134, 2.636121637993099, 2.691122513270625, 2.955269841410994, 2.877301803377398, 2.834259603310244, 2.729601826822466, 3.047539096860289, 3.082692313496543)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.09424102661766609, 0.09977962613732849, 0.1082786504442472, 0.09651939173358826, 0.09368764433342669, 0.09551283824153911, 0.1055046086301767, 0.09881333905530942, 0.09355860355896017, 0.09449432935074675, 0.09912181679367249, 0.09545310558251052, 0.09404438499955003, 0.09335271613782004, 0.09885468088260096, 0.09758783257704784, 0.09146391809895205, 0.09257540792357455, 0.09770720923873499, 0.09622654019154742, 0.09539720402200169, 0.09334387805636912, 0.0994247135581523, 0.1000694699481013)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.1435187710999097, 0.1367435091130798, 0.1383740820319155, 0.1392811574997546, 0.1426734207265342, 0.1373326219054835, 0.1392425432381202, 0.1376778939599666, 0.14257954701056, 0.1364559596819848, 0.1400995775036807, 0.1379705168934394, 0.1425999645822325, 0.137915266962632, 0.1408069520162238, 0.1386440898792071, 0.1434974795546587, 0.1379862469270051, 0.1413452495010717, 0.1386677449945413, 0.1433047026291803, 0.1378263793381722, 0.1409582909905325, 0.1386466084822203)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.6566459975612584, 0.7296845516434415, 0.7825067299761608, 0.6929824067103858, 0.6566580085929262, 0.6954854346789794, 0.7577038322960834, 0.7177139060831503, 0.656185305119751, 0.6924895737127856, 0.7075097481365948, 0.6918369788831993, 0.6594979548211449, 0.6768845697345013, 0.7020582397892607, 0.7038730079448083, 0.6373904153773873, 0.6709031514752849, 0.6912663112741838, 0.6939359992861743, 0.6656948604740102, 0.6772569845090368, 0.7053484605941354, 0.7217592340957547)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(-0.2205268161733829, -0.2102866718221291, -0.2085689693108573, -0.2394530840661691, -0.2167801529330208, -0.2245450734895393, -0.2042234949728717, -0.2188131470265026, -0.2141310778676606, -0.213615500526781, -0.2134256827950419, -0.219584191619921, -0.203654971182071, -0.2161559442698382, -0.2097528967073956, -0.212584891538512, -0.2148166530819151, -0.2193642417102586, -0.2109005717946058, -0.2153645488893773, -0.2026713111622997, -0.2181510351441237, -0.2053268057918707, -0.204366479447986)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value2 = $rowVals[$i] }

$ws = $wb.Worksheets.Item("cs_idio_reversal_lo")
$rowVals = @(2.87006099708343, 3.495152666754187, 3.472549170703409, 4.359618496073874, 3.101889834031375, 3.696266698329722, 3.408657934510352, 3.859518833829827, 3.08463146361521, 3.051908973368888, 3.031502245573072, 3.55457259635737, 3.035364608864892, 3.045126556856324, 3.13429683156661, 3.341754012524349, 3.065679104727256, 2.914672003044656, 2.92328773113455, 3.417723895687611, 3.077244090016542, 2.95554672410954, 3.111470672950025, 3.390962598477486)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.09608762685864813, 0.1072708344873139, 0.1068924677069305, 0.1205540387355168, 0.1004193962417126, 0.1105613614787311, 0.1058132516723109, 0.1131371988199321, 0.1001048857226809, 0.09950514634883989, 0.09912884087655227, 0.1082570870800561, 0.09920019990719697, 0.09938027319173104, 0.1010066515168857, 0.1046674058092616, 0.09975807421455452, 0.09693965716270192, 0.09710316636883221, 0.1059672731569212, 0.09996988203479429, 0.09771241878653525, 0.1005934620248359, 0.1055117734033224)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.1447782156902814, 0.1413753415066154, 0.1408011409124799, 0.1410700116567065, 0.1442830467401749, 0.1404363173064871, 0.1408519991137935, 0.1421401401601836, 0.1441165784228165, 0.1411562540332406, 0.1413431838209807, 0.1436598328487202, 0.1438857579422203, 0.1405825620225231, 0.1408953410295164, 0.1445758920019947, 0.1442842734482882, 0.140529398490011, 0.1413398040285667, 0.1436916762444449, 0.1443113879024253, 0.1405099048913242, 0.1409855900633113, 0.1431842959746952)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.6636884313051958, 0.7587662271520973, 0.7591733065101615, 0.8545688578298608, 0.6959888809566652, 0.7872704411455255, 0.7512371307334094, 0.7959553064492065, 0.6946104800586382, 0.7049290662346964, 0.7013344272908459, 0.7535654534281364, 0.6894372405296182, 0.7069174993112521, 0.7168913519697265, 0.7239616810236768, 0.6913994978829624, 0.689817633920866, 0.6870192514856348, 0.7374628505039651, 0.6927373056822648, 0.695413030576811, 0.7135017272308694, 0.7368948716412963)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(-0.2189320907800016, -0.2333284109823578, -0.2228212820497718, -0.2146796445856017, -0.2064115700177886, -0.2150573456833039, -0.2354148694862036, -0.2286562863059066, -0.2338224760860355, -0.2391837412671656, -0.2291149439201294, -0.2160076056045639, -0.2214088489146118, -0.225731460565701, -0.1909509505200376, -0.222109969169465, -0.2352277666698426, -0.2345359201788111, -0.2305943014381535, -0.2128789968325112, -0.2215741102876931, -0.2399926666497118, -0.2323043337480686, -0.2190413467806984)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value2 = $rowVals[$i] }

$ws = $wb.Worksheets.Item("cs_sr_lo")
$rowVals = @(6.441042491761898, 7.01395702258098, 7.307037533940925, 7.185775508617017, 6.270950738503766, 7.039043749440713, 6.688119471090064, 7.293094645424436, 6.081957999730171, 6.258981766505527, 6.633726913895706, 6.676742892464875, 6.182766420795312, 6.330195229125548, 6.265927720029803, 6.430739403042142, 6.13367555777453, 6.047960063822844, 6.501500901022705, 6.500986570377511, 5.89639608072023, 5.996250966212785, 6.132093394225115, 6.081586468475341)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.1457604208483156, 0.1515366460204179, 0.1543442227864777, 0.1531939661480632, 0.1439655942291871, 0.1517806797631642, 0.1483006260331849, 0.1542127639180444, 0.1419248265426989, 0.1438378268352587, 0.147748015379656, 0.1481853457979712, 0.1430195999632942, 0.1445951494406228, 0.1439119979799852, 0.1456527951717836, 0.1424882777584831, 0.1415523333046833, 0.1463891816928651, 0.1463838526627574, 0.1398710982148514, 0.1409825667244624, 0.1424710969873204, 0.1419207649336114)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.1257806906015005, 0.131576834112743, 0.1353477106980147, 0.1390898650583813, 0.1227079564505745, 0.1281562516920592, 0.1313167985153713, 0.1351540826821936, 0.1261179275542171, 0.1315184693857204, 0.1334549267938254, 0.136583965670805, 0.125345484796746, 0.1301582450137256, 0.1325364707720482, 0.1349910164938986, 0.1262400316460429, 0.1315850789947206, 0.1333684975161618, 0.1360119762375059, 0.1261920697359728, 0.1309880595705536, 0.1328369526693273, 0.1349251452176672)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(1.15884576679671, 1.151697006864994, 1.140353405244125, 1.10140279511927, 1.173237648099657, 1.184340816457952, 1.129334766837357, 1.141014469245935, 1.125334274793617, 1.093670170486913, 1.107100493996089, 1.084939546675104, 1.141003205621707, 1.110918093781725, 1.085829410890999, 1.078981394131268, 1.128709141629477, 1.075747602890162, 1.097629383394122, 1.076257082002397, 1.108398479456742, 1.076300902438557, 1.07252608648721, 1.051848154060969)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(-0.1459081367157734, -0.1576794365059824, -0.1704721054351646, -0.2475755090726953, -0.1643194727445826, -0.1604919876312526, -0.1926049935635472, -0.217971850232228, -0.1591079456626258, -0.1630078801914067, -0.1794903434813099, -0.2319170522804602, -0.1727766327656879, -0.1771728208391011, -0.1929049118812454, -0.2236336886292617, -0.1683878981125063, -0.1727636439395525, -0.1841260568188997, -0.2377729714390679, -0.1847341010246271, -0.1849293544509386, -0.1978407989534848, -0.2283425346471357)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value2 = $rowVals[$i] }

$ws = $wb.Worksheets.Item("cs_idio_sr_lo")
$rowVals = @(2.795520683618175, 3.525298907899193, 4.212862881528202, 5.712220159192655, 2.78120303660786, 3.539179294087547, 3.992007218603952, 5.262812946031159, 2.998299784715079, 3.501466564952408, 3.717996218191871, 4.847656264130956, 2.940759456603477, 3.455771867937685, 3.694948479732306, 4.622416694692082, 3.179906902845833, 3.598778044007302, 3.69275053887023, 4.698229710929669, 3.028748028026499, 3.537040226762368, 3.687731251672846, 4.54738198892745)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.09464332905775286, 0.1077727108559423, 0.1184468287996547, 0.1377811252330596, 0.09436288535762105, 0.1080027452660404, 0.1151690093624436, 0.1324479949170241, 0.09851276706538825, 0.1073762078420268, 0.1109089881387124, 0.1271942599853013, 0.09743371948991286, 0.1066104730135349, 0.1105402245419917, 0.1241965269844179, 0.101825934429483, 0.1089830589076033, 0.1105049694193947, 0.1252178343306887, 0.09907791643808528, 0.1079673380856938, 0.1104244019597902, 0.1231729846609035)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.1235583349772551, 0.1342162091472668, 0.1343440908608251, 0.1348860518464058, 0.1231284787435215, 0.1329436946948639, 0.1353656293313746, 0.1358878870161987, 0.1279554307630029, 0.1357217069346307, 0.134746129619399, 0.1338681282961094, 0.1283131331545636, 0.1354615032359929, 0.1353172226205906, 0.1352936701778882, 0.129474097386966, 0.1362900264097181, 0.1359235926932003, 0.134707933636649, 0.1293983810103037, 0.1357111518753726, 0.1354960726169553, 0.1353689861872033)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(0.7659809358484397, 0.8029783551529922, 0.8816675749613782, 1.021463104205544, 0.7663774158550303, 0.8123946420620499, 0.8507994971198356, 0.9746858077293923, 0.7698990693708976, 0.7911498482239379, 0.8230959097080083, 0.9501459503785257, 0.7593433119004743, 0.7870167572834662, 0.8168969359645378, 0.9179773659855667, 0.7864579594260502, 0.7996407497932092, 0.8129932944666994, 0.9295505539298213, 0.7656812679147506, 0.795567177742646, 0.8149638570850518, 0.9099054970432168)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value2 = $rowVals[$i] }
$rowVals = @(-0.1758129042317258, -0.1729465975100464, -0.1785164495277428, -0.1711611834787135, -0.1718250645359084, -0.177762293580551, -0.1959650665262684, -0.1778933617689102, -0.1871587146284003, -0.1830244938601768, -0.1806675625442244, -0.1806675625442243, -0.1835914019524187, -0.1791695072091938, -0.1933510182497865, -0.1822000334986038, -0.1826646269758735, -0.1801292616216639, -0.1825220071704523, -0.1825220071704523, -0.182165213731955, -0.1797512825505978, -0.1814141106269513, -0.1806415977544142)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value2 = $rowVals[$i] }
